# Questions_a_traiter.xlsx — fill in the "Réponse" column for three rows
# (PMD / Maven et Eclipse / Exceptions) that were answered after the Maven
# components were added, per commit "COMMIT_DL_DOMICILE après ajout
# composants Maven".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Row 8 — Exceptions: strategy for typing/propagating exceptions to the view.
$ws.Range("E8").Value = "Typer les Exceptions`nToute la logique métier throw les Exceptions jusqu' aux contrôleurs chargés de transmettre à la vue (messages des exceptions humanisés pour les utilisateurs)"

# Row 3 — PMD: how to export/import the PMD ruleset.
$ws.Range("E3").Value = "Exporter ou importer : \Workspace\.metadata\.plugins\net.sourceforge.pmd.eclipse.plugin\ruleset.xml "

# Row 4 — Maven et Eclipse: the maven-compiler-plugin snippet that pins the JDK.
$ws.Range("E4").Value = "<plugin>`n   <groupId>org.apache.maven.plugins</groupId>`n   <artifactId>maven-compiler-plugin</artifactId>`n   <configuration>`n      <source>1.8</source>`n      <target>1.8</target>`n   </configuration>`n</plugin> `n"

# Row heights grew to fit the newly-added answers.
$ws.Rows(4).RowHeight = 110.25
$ws.Rows(8).RowHeight = 53.25

# Scroll the view up a bit and leave the selection on E7, matching where the
# author ended up after entering the last answer.
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E7").Select()
